$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G = "Recorded By" holds a comma-separated list of recorder names/
# emails (e.g. "dnasr281@gmail.com, System"). This pass rotates each such
# list left by one position (first entry moves to the end), e.g.
#   "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#   "backup@backdoor.com, System, system" -> "System, system, backup@backdoor.com"
# Cells holding only a single value (no comma) are left untouched.

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
            $newVal = $rotated -join ", "
            $cell.Value = $newVal
        }
    }
}
